{"js": "// Replace the 25 multiplication-problem cell texts (5 per row, across the\n// 5 \"problem\" rows of the single table) with their new values, in document\n// order. Blank spacer rows are left untouched. Cell run formatting\n// (font/size) is preserved because we only change the cell's text value,\n// not its formatting.\nconst replacements = [\n  [\"402\u00d76=\", \"704\u00d75=\"],\n  [\"914\u00d76=\", \"220\u00d73=\"],\n  [\"410\u00d76=\", \"886\u00d75=\"],\n  [\"199\u00d74=\", \"110\u00d78=\"],\n  [\"940\u00d73=\", \"687\u00d73=\"],\n  [\"549\u00d75=\", \"942\u00d72=\"],\n  [\"275\u00d74=\", \"638\u00d77=\"],\n  [\"461\u00d79=\", \"123\u00d76=\"],\n  [\"665\u00d76=\", \"809\u00d77=\"],\n  [\"853\u00d72=\", \"901\u00d76=\"],\n  [\"330\u00d75=\", \"687\u00d73=\"],\n  [\"580\u00d77=\", \"991\u00d79=\"],\n  [\"294\u00d78=\", \"484\u00d72=\"],\n  [\"408\u00d76=\", \"624\u00d74=\"],\n  [\"484\u00d78=\", \"133\u00d74=\"],\n  [\"846\u00d76=\", \"216\u00d78=\"],\n  [\"995\u00d77=\", \"278\u00d75=\"],\n  [\"695\u00d78=\", \"246\u00d76=\"],\n  [\"461\u00d77=\", \"480\u00d72=\"],\n  [\"342\u00d77=\", \"731\u00d78=\"],\n  [\"587\u00d78=\", \"238\u00d72=\"],\n  [\"705\u00d73=\", \"522\u00d76=\"],\n  [\"395\u00d78=\", \"559\u00d75=\"],\n  [\"987\u00d78=\", \"963\u00d72=\"],\n  [\"214\u00d79=\", \"686\u00d75=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Gather every cell across every row (in order), then load each cell's\n// current text value.\nconst cells = [];\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nfor (const row of rows.items) {\n  for (const cell of row.cells.items) {\n    cell.load(\"value\");\n    cells.push(cell);\n  }\n}\nawait context.sync();\n\n// Walk the cells in document order and apply the next replacement whenever\n// the cell's current text matches the expected \"old\" value.\nlet idx = 0;\nfor (const cell of cells) {\n  if (idx < replacements.length && cell.value === replacements[idx][0]) {\n    cell.value = replacements[idx][1];\n    idx++;\n  }\n}\nawait context.sync();\n", "ps1": "# Replace the 25 multiplication-problem cell texts (5 per row, across the\n# 5 \"problem\" rows of the single table) with their new values, in document\n# order. Blank spacer rows are left untouched. Only the cell's Range.Text\n# is changed, so the existing run formatting (font/size) is preserved.\n\n$replacements = @(\n    @(\"402\u00d76=\", \"704\u00d75=\"),\n    @(\"914\u00d76=\", \"220\u00d73=\"),\n    @(\"410\u00d76=\", \"886\u00d75=\"),\n    @(\"199\u00d74=\", \"110\u00d78=\"),\n    @(\"940\u00d73=\", \"687\u00d73=\"),\n    @(\"549\u00d75=\", \"942\u00d72=\"),\n    @(\"275\u00d74=\", \"638\u00d77=\"),\n    @(\"461\u00d79=\", \"123\u00d76=\"),\n    @(\"665\u00d76=\", \"809\u00d77=\"),\n    @(\"853\u00d72=\", \"901\u00d76=\"),\n    @(\"330\u00d75=\", \"687\u00d73=\"),\n    @(\"580\u00d77=\", \"991\u00d79=\"),\n    @(\"294\u00d78=\", \"484\u00d72=\"),\n    @(\"408\u00d76=\", \"624\u00d74=\"),\n    @(\"484\u00d78=\", \"133\u00d74=\"),\n    @(\"846\u00d76=\", \"216\u00d78=\"),\n    @(\"995\u00d77=\", \"278\u00d75=\"),\n    @(\"695\u00d78=\", \"246\u00d76=\"),\n    @(\"461\u00d77=\", \"480\u00d72=\"),\n    @(\"342\u00d77=\", \"731\u00d78=\"),\n    @(\"587\u00d78=\", \"238\u00d72=\"),\n    @(\"705\u00d73=\", \"522\u00d76=\"),\n    @(\"395\u00d78=\", \"559\u00d75=\"),\n    @(\"987\u00d78=\", \"963\u00d72=\"),\n    @(\"214\u00d79=\", \"686\u00d75=\")\n)\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$idx = 0\nfor ($r = 1; $r -le $tbl.Rows.Count; $r++) {\n    for ($c = 1; $c -le $tbl.Columns.Count; $c++) {\n        if ($idx -ge $replacements.Length) { break }\n\n        $cell = $tbl.Cell($r, $c)\n        $raw = $cell.Range.Text\n        # Strip the trailing end-of-cell marker (CR + BEL) before comparing.\n        $text = $raw -replace \"`r`a$\", \"\"\n\n        $expectedOld = $replacements[$idx][0]\n        if ($text -eq $expectedOld) {\n            $cell.Range.Text = $replacements[$idx][1]\n            $idx++\n        }\n    }\n}\n"}
